# Auto-generated edit script applying market-data refresh from diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H33").Value = 355.375
$ws.Range("I33").Value = 355.375
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 355.375
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -126.375

$ws.Range("H70").Value = 1491.1666
$ws.Range("I70").Value = 950
$ws.Range("J70").Value = 1599.4
$ws.Range("K70").Value = 2850
$ws.Range("L70").Value = 4798.200000000001
$ws.Range("M70").Value = -2580
$ws.Range("N70").Value = -5338.200000000001

$ws.Range("H73").Value = 1491.1666
$ws.Range("I73").Value = 950
$ws.Range("J73").Value = 1599.4
$ws.Range("K73").Value = 2850
$ws.Range("L73").Value = 4798.200000000001
$ws.Range("M73").Value = -1914
$ws.Range("N73").Value = -6670.200000000001

$ws.Range("H94").Value = 3768.3333
$ws.Range("I94").Value = 3768.3333
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 3768.3333
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -3317.3333

$ws.Range("H96").Value = 613
$ws.Range("I96").Value = 613
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1839
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -466
$ws.Range("N96").ClearContents()

$ws.Range("H98").Value = 738.8
$ws.Range("I98").Value = 686.25
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 686.25
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 811.75

$ws.Range("H117").Value = 90027.55
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 90027.55
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 90027.55
$ws.Range("N117").Value = -99205.55

$ws.Range("H122").Value = 738.8
$ws.Range("I122").Value = 686.25
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2058.75
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = 391.25

$ws.Range("H132").Value = 2195.8438
$ws.Range("I132").Value = 1675.7
$ws.Range("J132").Value = 9998
$ws.Range("K132").Value = 5027.1
$ws.Range("L132").Value = 29994
$ws.Range("M132").Value = -2497.1

$ws.Range("H137").Value = 662508.0600000001
$ws.Range("I137").Value = 2648.5833
$ws.Range("J137").Value = 1454339.4
$ws.Range("K137").Value = 7945.749899999999
$ws.Range("L137").Value = 4363018.199999999
$ws.Range("M137").Value = -5395.749899999999
$ws.Range("N137").Value = -4368118.199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7612.3403
$ws.Range("I32").Value = 3679.6758
$ws.Range("J32").Value = 22163.2
$ws.Range("K32").Value = 3679.6758
$ws.Range("L32").Value = 22163.2
$ws.Range("M32").Value = -3392.6758

$ws.Range("H97").Value = 699.875
$ws.Range("I97").Value = 699.875
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 699.875
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -203.875

$ws.Range("H102").Value = 46421.08
$ws.Range("I102").Value = 49750.81
$ws.Range("J102").Value = 28940
$ws.Range("K102").Value = 49750.81
$ws.Range("L102").Value = 28940
$ws.Range("M102").Value = -48128.81
$ws.Range("N102").Value = -32184

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 61731.35
$ws.Range("I22").Value = 74490.07000000001
$ws.Range("J22").Value = 2190.6667
$ws.Range("K22").Value = 74490.07000000001
$ws.Range("L22").Value = 2190.6667
$ws.Range("M22").Value = -74317.07000000001
$ws.Range("N22").Value = -2536.6667

$ws.Range("H94").Value = 1876.5769
$ws.Range("I94").Value = 1850.409
$ws.Range("J94").Value = 2020.5
$ws.Range("K94").Value = 1850.409
$ws.Range("L94").Value = 2020.5
$ws.Range("M94").Value = -1399.409
$ws.Range("N94").Value = -2922.5

$ws.Range("H102").Value = 44000
$ws.Range("I102").Value = 44000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 44000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -40755

$ws.Range("H105").Value = 94886.27
$ws.Range("I105").Value = 127874.875
$ws.Range("J105").Value = 6916.6665
$ws.Range("K105").Value = 127874.875
$ws.Range("L105").Value = 6916.6665
$ws.Range("M105").Value = -126127.875

$ws.Range("H134").Value = 6161.8423
$ws.Range("I134").Value = 3890.3845
$ws.Range("J134").Value = 11083.333
$ws.Range("K134").Value = 11671.1535
$ws.Range("L134").Value = 33249.999
$ws.Range("M134").Value = -9136.1535

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8549361
$ws.Range("I99").Value = 12347912
$ws.Range("J99").Value = 2622
$ws.Range("K99").Value = 12347912
$ws.Range("L99").Value = 2622
$ws.Range("M99").Value = -12346414
$ws.Range("N99").Value = -5618

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H126").Value = 8549361
$ws.Range("I126").Value = 12347912
$ws.Range("J126").Value = 2622
$ws.Range("K126").Value = 37043736
$ws.Range("L126").Value = 7866
$ws.Range("M126").Value = -37041266
$ws.Range("N126").Value = -12806

$ws.Range("H132").Value = 1569647.8
$ws.Range("I132").Value = 1685677.2
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 5057031.6
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -5054501.6
$ws.Range("N132").Value = -14810

$ws.Range("H134").Value = 1686579.8
$ws.Range("I134").Value = 2103035.5
$ws.Range("J134").Value = 113302.336
$ws.Range("K134").Value = 6309106.5
$ws.Range("L134").Value = 339907.008
$ws.Range("M134").Value = -6306571.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 2033198.1
$ws.Range("I9").Value = 3343336.2
$ws.Range("J9").Value = 1050594.5
$ws.Range("K9").Value = 10030008.6
$ws.Range("L9").Value = 3151783.5
$ws.Range("M9").Value = -10029784.6
$ws.Range("N9").Value = -3152231.5

$ws.Range("H18").Value = 10818.6
$ws.Range("I18").Value = 11798.444
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 35395.33199999999
$ws.Range("L18").Value = 6000
$ws.Range("M18").Value = -35226.33199999999

$ws.Range("H26").Value = 17705.5
$ws.Range("I26").Value = 510
$ws.Range("J26").Value = 26303.25
$ws.Range("K26").Value = 1530
$ws.Range("L26").Value = 78909.75
$ws.Range("M26").Value = -1242
$ws.Range("N26").Value = -79485.75

$ws.Range("H57").Value = 4499.1665
$ws.Range("I57").Value = 3995
$ws.Range("J57").Value = 4600
$ws.Range("K57").Value = 11985
$ws.Range("L57").Value = 13800
$ws.Range("M57").Value = -11426
$ws.Range("N57").Value = -14918

$ws.Range("H97").Value = 112.46667
$ws.Range("I97").Value = 110.125
$ws.Range("J97").Value = 115.14286
$ws.Range("K97").Value = 330.375
$ws.Range("L97").Value = 345.42858
$ws.Range("M97").Value = 165.625
$ws.Range("N97").Value = -1337.42858

$ws.Range("H128").Value = 378275
$ws.Range("I128").Value = 378275
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 1134825
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -1129845

$ws.Range("H132").Value = 3891.1365
$ws.Range("I132").Value = 991.5714
$ws.Range("J132").Value = 5244.2666
$ws.Range("K132").Value = 8924.142600000001
$ws.Range("L132").Value = 47198.3994
$ws.Range("M132").Value = -6394.142600000001
$ws.Range("N132").Value = -52258.3994

$ws.Range("H137").Value = 3853.0527
$ws.Range("I137").Value = 2131.25
$ws.Range("J137").Value = 6804.7144
$ws.Range("K137").Value = 6393.75
$ws.Range("L137").Value = 20414.1432
$ws.Range("M137").Value = -1293.75
$ws.Range("N137").Value = -30614.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 6318.4
$ws.Range("I18").Value = 3866.3333
$ws.Range("J18").Value = 9996.5
$ws.Range("K18").Value = 3866.3333
$ws.Range("L18").Value = 9996.5
$ws.Range("M18").Value = -3573.3333

$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()

$ws.Range("H70").Value = 9497.5
$ws.Range("I70").Value = 8746.25
$ws.Range("J70").Value = 11000
$ws.Range("K70").Value = 8746.25
$ws.Range("L70").Value = 11000
$ws.Range("M70").Value = -8476.25

$ws.Range("H73").Value = 9497.5
$ws.Range("I73").Value = 8746.25
$ws.Range("J73").Value = 11000
$ws.Range("K73").Value = 8746.25
$ws.Range("L73").Value = 11000
$ws.Range("M73").Value = -7810.25

$ws.Range("H122").Value = 4286.8887
$ws.Range("I122").Value = 5020.5
$ws.Range("J122").Value = 3700
$ws.Range("K122").Value = 15061.5
$ws.Range("L122").Value = 11100
$ws.Range("M122").Value = -12611.5
$ws.Range("N122").Value = -16000

$ws.Range("H132").Value = 5743.625
$ws.Range("I132").Value = 2349.5
$ws.Range("J132").Value = 6875
$ws.Range("K132").Value = 7048.5
$ws.Range("L132").Value = 20625
$ws.Range("M132").Value = -4518.5
$ws.Range("N132").Value = -25685

$ws.Range("H135").Value = 70000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 70000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

$ws.Range("H136").Value = 83916.55499999999
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 83916.55499999999
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 251749.665
$ws.Range("N136").Value = -256849.665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2678.4285
$ws.Range("I22").Value = 1875
$ws.Range("J22").Value = 2999.8
$ws.Range("K22").Value = 1875
$ws.Range("L22").Value = 2999.8
$ws.Range("M22").Value = -1580
$ws.Range("N22").Value = -3589.8

$ws.Range("H27").Value = 2678.4285
$ws.Range("I27").Value = 1875
$ws.Range("J27").Value = 2999.8
$ws.Range("K27").Value = 1875
$ws.Range("L27").Value = 2999.8
$ws.Range("M27").Value = -1768
$ws.Range("N27").Value = -3213.8

$ws.Range("H96").Value = 62112.5
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 62112.5
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 62112.5
$ws.Range("N96").Value = -67604.5

$ws.Range("H117").Value = 89095.336
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 89095.336
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 89095.336
$ws.Range("N117").Value = -98273.336
$ws.Range("M117").ClearContents()

$ws.Range("H132").Value = 2351
$ws.Range("I132").Value = 2363.3635
$ws.Range("J132").Value = 2215
$ws.Range("K132").Value = 7090.0905
$ws.Range("L132").Value = 6645
$ws.Range("M132").Value = -4560.0905

$ws.Range("H136").Value = 4196.857
$ws.Range("I136").Value = 5470.875
$ws.Range("J136").Value = 2498.1667
$ws.Range("K136").Value = 16412.625
$ws.Range("L136").Value = 7494.500100000001
$ws.Range("M136").Value = -13862.625
$ws.Range("N136").Value = -12594.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H136").Value = 1506.6
$ws.Range("I136").Value = 1227.6666
$ws.Range("J136").Value = 1925
$ws.Range("K136").Value = 3682.9998
$ws.Range("L136").Value = 5775
$ws.Range("M136").Value = -1132.9998
$ws.Range("N136").Value = -10875
